$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

$ws.Range("A54").Value = "agregar campo peso a la tabla de articulos"
$ws.Range("B54").Value = "no comenzado"

$ws.Range("A50").Select()
